# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet between the
# existing "Late" column (N) and the "Outstanding" columns (old O/P), pushing
# the old N/O/P data one column to the right (-> O/P/Q), then leave the view
# focused on the repayment schedule sheet the way the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N; existing N/O/P columns (and their data/styles)
# shift right to O/P/Q automatically.
[void]$ws.Columns("N").Insert()

# The inserted column picks up a width close to its neighbouring "M" column
# (same as what Excel does when a column is inserted next to formatted data).
$ws.Columns("N").ColumnWidth = 9.85

# Leave the selection where the author left it after making the edit, which
# also marks "Repayment schedule" as the active sheet/tab.
[void]$ws.Range("R8").Select()
